$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '66.049.44'
$c.NumberFormat = 'General'
$ws.Range('E2').Value = '  -4.92%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.339.26'
$c.NumberFormat = 'General'
$ws.Range('E3').Value = '  -5.75%  '
$ws.Range('E4').Value = '  +0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '559.93'
$c.NumberFormat = 'General'
$ws.Range('E5').Value = '  -4.38%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '181.99'
$c.NumberFormat = 'General'
$ws.Range('E6').Value = '  -7.96%  '
$ws.Range('E7').Value = '  +0.01%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.591'
$c.NumberFormat = 'General'
$ws.Range('E8').Value = '  -3.46%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '3.329.50'
$c.NumberFormat = 'General'
$ws.Range('E9').Value = '  -5.63%  '
$ws.Range('E10').Value = '  -9.69%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.587'
$c.NumberFormat = 'General'
$ws.Range('E11').Value = '  -6.88%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '47.37'
$c.NumberFormat = 'General'
$ws.Range('E12').Value = '  -8.81%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000264'
$c.NumberFormat = 'General'
$ws.Range('E13').Value = '  -7.90%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '8.63'
$c.NumberFormat = 'General'
$ws.Range('E14').Value = '  -6.82%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.875.44'
$c.NumberFormat = 'General'
$ws.Range('E15').Value = '  -5.60%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '604.87'
$c.NumberFormat = 'General'
$ws.Range('E16').Value = '  -8.84%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '66.131.44'
$c.NumberFormat = 'General'
$ws.Range('E17').Value = '  -4.96%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '18.04'
$c.NumberFormat = 'General'
$ws.Range('E18').Value = '  -2.69%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '3.343.14'
$c.NumberFormat = 'General'
$ws.Range('E19').Value = '  -5.54%  '
$ws.Range('E20').Value = '  -3.88%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '11.43'
$c.NumberFormat = 'General'
$ws.Range('E21').Value = '  -8.14%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.906'
$c.NumberFormat = 'General'
$ws.Range('E22').Value = '  -6.47%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '16.87'
$c.NumberFormat = 'General'
$ws.Range('E23').Value = '  -7.78%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.01'
$c.NumberFormat = 'General'
$ws.Range('E24').Value = '  -5.05%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '100.18'
$c.NumberFormat = 'General'
$ws.Range('E25').Value = '  -4.97%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '4.05'
$c.NumberFormat = 'General'
$ws.Range('E26').Value = '  -7.53%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '6.00'
$c.NumberFormat = 'General'
$ws.Range('E27').Value = '  -0.98%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.67'
$c.NumberFormat = 'General'
$ws.Range('E28').Value = '  -8.35%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.29'
$c.NumberFormat = 'General'
$ws.Range('E29').Value = '  -8.84%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '8.70'
$c.NumberFormat = 'General'
$ws.Range('E30').Value = '  -9.92%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '30.35'
$c.NumberFormat = 'General'
$ws.Range('E31').Value = '  -9.29%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '6.25'
$c.NumberFormat = 'General'
$ws.Range('E32').Value = '  -8.44%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.77'
$c.NumberFormat = 'General'
$ws.Range('E33').Value = '  -15.20%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '11.03'
$c.NumberFormat = 'General'
$ws.Range('E34').Value = '  -6.82%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.843.14'
$c.NumberFormat = 'General'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.105'
$c.NumberFormat = 'General'
$ws.Range('E36').Value = '  -5.71%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '537.17'
$c.NumberFormat = 'General'
$ws.Range('E37').Value = '  +6.52%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '57.56'
$c.NumberFormat = 'General'
$ws.Range('E38').Value = '  -6.98%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.NumberFormat = 'General'
$ws.Range('E39').Value = '  -0.08%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.39'
$c.NumberFormat = 'General'
$ws.Range('E40').Value = '  -8.77%  '
$ws.Range('D41').Value = '0.0₃0711'
$ws.Range('E41').Value = '  -12.74%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.65'
$c.NumberFormat = 'General'
$ws.Range('E42').Value = '  -9.42%  '
$ws.Range('E43').Value = '  -7.06%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.342'
$c.NumberFormat = 'General'
$ws.Range('E44').Value = '  -8.49%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '32.00'
$c.NumberFormat = 'General'
$ws.Range('E45').Value = '  -7.68%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.15'
$c.NumberFormat = 'General'
$ws.Range('E46').Value = '  +16.89%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0411'
$c.NumberFormat = 'General'
$ws.Range('E47').Value = '  -9.05%  '
$ws.Range('E48').Value = '  -7.86%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.129'
$c.NumberFormat = 'General'
$ws.Range('E49').Value = '  -5.00%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.61'
$c.NumberFormat = 'General'
$ws.Range('E50').Value = '  -9.28%  '
$ws.Range('E51').Value = '  -0.02%  '
